{"js": "const paras = context.document.body.paragraphs;\nparas.load(\"text\");\nawait context.sync();\nconst p = paras.items[78];\nconst rng = p.getRange(\"Start\");\nconst xml = '<?xml version=\"1.0\" standalone=\"yes\"?><pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\"><pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\"><pkg:xmlData><w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\"><w:body><w:p><w:r><w:lastRenderedPageBreak/></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>';\nrng.insertOoxml(xml, \"Start\");\nawait context.sync();\nreturn \"done\";\n", "ps1": "$d = $word.ActiveDocument\n$full = $d.WordOpenXML\n$idx = $full.IndexOf(\"_Toc145419454\")\nWrite-Output $full.Substring($idx-100, 250)\n"}
